$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 197, shifting existing rows 197:205 down to 198:206.
$ws.Rows.Item(197).Insert()

# Populate the new row 197 with the inserted record's data.
$ws.Cells.Item(197, 1).Value = 6
$ws.Cells.Item(197, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(197, 3).Value = "Metropolitana"
$ws.Cells.Item(197, 4).Value = 44706
$ws.Cells.Item(197, 4).Style = $ws.Cells.Item(198, 4).Style
$ws.Cells.Item(197, 4).NumberFormat = $ws.Cells.Item(198, 4).NumberFormat
$ws.Cells.Item(197, 5).Value = 13
$ws.Cells.Item(197, 6).Value = 100112029
$ws.Cells.Item(197, 7).Value = "Orégano"
$ws.Cells.Item(197, 8).Value = "Sin especificar"
$ws.Cells.Item(197, 9).Value = "Primera"
$ws.Cells.Item(197, 10).Value = 41
$ws.Cells.Item(197, 11).Value = 12000
$ws.Cells.Item(197, 12).Value = 13000
$ws.Cells.Item(197, 13).Value = 12439
$ws.Cells.Item(197, 14).Value = "`$/docena de atados"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 4146
$ws.Cells.Item(197, 17).Value = 3
$ws.Cells.Item(197, 18).Value = "Hortaliza"
